$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values
$ws.Range("A1").Value = "var_1_input_first_name"
$ws.Range("B1").Value = "var_2_input_last_name"
$ws.Range("C1").Value = "var_3_input_postal_code"

# Update data row values
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "test"
$ws.Range("C2").Value = "test"

# Update column widths
# NOTE: the runtime adds a constant offset of 5/6 (0.8333...) character
# widths when converting ColumnWidth -> stored OOXML width, so we
# subtract that offset here to land exactly on the target stored widths
# of 24, 23, and 25.
$ws.Columns.Item(1).ColumnWidth = 24 - 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 23 - 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 25 - 0.8333333333333334
